# Updates cryptos list values (Price and Volume(1h) columns, and a few Coin/Link swaps)
# per commit: "Updated cryptos list on Sun Dec 10 11:36:57 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.734.18'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.337.12'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.78'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.664'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.78'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -6.71%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.596'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0989'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.86'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.15'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.97%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.12'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -6.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.686.63'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.02'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.895'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.343.03'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.681.72'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '77.76'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.96'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.03%  '
$ws.Range("B24").Value = 'ImmutableX'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.90'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.65%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.71'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.65%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.28'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -7.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '176.43'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.11'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.07%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.132'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.76%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.32%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.31'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.77%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.79%  '
$ws.Range("B38").Value = 'FTXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.88'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +30.56%  '
$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.35'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.35'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.45%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '66.51'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +19.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.13'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.62'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.69%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.87%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.21'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.40'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.58%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.59'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.19%  '
